# Update the "N=20(Solution)" worksheet with refreshed measurement data
# (columns C and D, rows 3-22) as captured from a later simulation run,
# then restore each sheet's selection / active-tab view state.

$wb = $excel.ActiveWorkbook
$wsOriginal = $wb.Worksheets.Item("N=20(Original)")
$wsSolution = $wb.Worksheets.Item("N=20(Solution)")
$wsOptimal  = $wb.Worksheets.Item("N=20(Optimal)")

# --- Data update: N=20(Solution)!C3:D22 -------------------------------
$newData = @(
  @(404,376),
  @(405,416),
  @(411,375),
  @(425,386),
  @(412,355),
  @(404,407),
  @(400,401),
  @(430,402),
  @(392,424),
  @(424,369),
  @(379,406),
  @(384,399),
  @(399,387),
  @(416,390),
  @(375,381),
  @(387,418),
  @(429,428),
  @(384,397),
  @(417,430),
  @(420,416)
)

$rowCount = $newData.Count
$arr = New-Object 'object[,]' $rowCount,2
for ($i = 0; $i -lt $rowCount; $i++) {
    $arr[$i,0] = $newData[$i][0]
    $arr[$i,1] = $newData[$i][1]
}
$wsSolution.Range("C3:D22").Value = $arr

# --- View state: selections on each sheet ------------------------------
$wsOriginal.Activate() | Out-Null
$wsOriginal.Range("M21").Select() | Out-Null

$wsOptimal.Activate() | Out-Null
$wsOptimal.Range("P26").Select() | Out-Null

$wsSolution.Activate() | Out-Null
$wsSolution.Range("E42").Select() | Out-Null

$excel.Calculate() | Out-Null
